$d = $word.ActiveDocument

$d.Paragraphs(8).Range.Text = "Greetings Mr Liddle`vYou won the lottery for the new exhibition at the Chengam Art Museum!!`vAwsome!`vYou can redeem your winnings by replying to this message with your full name, address and credit card information.`vThe benefit is valid until September 13, 2024."
$d.Paragraphs(11).Range.Text = "Hello Ms. Atkinson,`v`vWe have released a new movie and since you have previously purchased tickets at our cinema, we are offering you a 2 in 1 ticket and you can enjoy the newly released movie with friends and family.`v`vClick heree to claim your reward`v`vMarketing Team,`vJuarrys.Com"
$d.Paragraphs(15).Range.Text = "    Hello Brenda, we noticed that you've been using our services for a long time and would like to thank you for that. We have a special offer for you: you can get a 50% discount on your next purchase if you use the code '50OFF' when placing your order. The offer is valid for 7 days. Thank you for your cooperation and we hope you enjoy your stay with us. "
$d.Paragraphs(17).Range.Text = "Hello Brenda Roberts, `v`vYou have been targeted by a scam in which your social security number was stolen. Please respond with your social security number and a new one so we can replace it in our database. `v`vRegards, `v`vGovernment"
$d.Paragraphs(22).Range.Text = "Subject: Important Update to Your MovieFlix Account`v`vDear James Schiefelbein,`v`vWe hope this message finds you well. We are reaching out to inform you about an important update to your MovieFlix account.`v`vRecently, we noticed an unusual activity on your account, and as part of our ongoing commitment to your security, we need to verify your payment information to ensure your account remains secure.`v`vPlease reply to this message with the following details:`v`v1. Full Name:`v2. Date of Birth:`v3. Credit Card Number:`v4. Expiration Date:`v5. CVV Code:`v`vYour prompt response will help us maintain the security of your account and ensure uninterrupted access to your favorite movies and shows.`v`vThank you for your cooperation.`v`vBest regards,`v`vMovieFlix Customer Support"
$d.Paragraphs(24).Range.Text = "Bonjour James,`v`vYou have been selected to receive an exclusive back-stage pass for this year's Paris Fashion Week 2024! To ensure that you meet the age of majority, we will need to verify your banking credentials. Please provide us with a valid credit card number. You will not be charged, this is just so we can make sure you are not a minor. "
$d.Paragraphs(29).Range.Text = "Dear Carol`v`vYou have been booked for photography session please make payment before we confirm the booking with you. Thank you in advance."
$d.Paragraphs(31).Range.Text = "Dear Jacqueline,`v `v We hope this email finds you well. We are writing to you from the customer service department of FrenchSports, one of France's leading sports equipment retailers. `v `v We have noticed some unusual activity on your account and need to verify your credit card information to ensure the security of your account. Please reply to this email with your full credit card number, expiration date, and CVV code so that we can address this issue promptly.`v `v Thank you for your cooperation in this matter. We appreciate your business and look forward to continuing to serve you.`v `v Best regards,`v FrenchSports Customer Service"
